$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Event Data")

# Test setting values
$ws.Range("T1").Value = "North/South/??"
$ws.Range("U1").Value = "Notes"
Write-Host ("T1 after set: " + $ws.Range("T1").Value())
Write-Host ("U1 after set: " + $ws.Range("U1").Value())

# Test setting font color / interior
$c = $ws.Range("C2")
$c.Font.Color = 255
Write-Host ("C2 font color: " + $c.Font.Color())

$c.Interior.Color = 65535
Write-Host ("C2 interior color: " + $c.Interior.Color())

# column width
$ws.Columns.Item(21).ColumnWidth = 25.83203125
Write-Host ("U width after: " + $ws.Range("U1").ColumnWidth())

# Entire row format
$ws.Rows.Item(7).Font.Size = 8
Write-Host ("Row7 A7 font size: " + $ws.Range("A7").Font.Size())
